$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (not a number) into a cell, even if the
# value looks numeric (e.g. "1"), without leaving the cell's number format
# changed afterwards.
function Set-TextValue($rangeRef, $text) {
    $rng = $ws.Range($rangeRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Insert 5 new rows before the current row 321 (STAT_PAZUCAR), shifting the
# existing rows 321-327 down to 326-332.
for ($i = 0; $i -lt 5; $i++) {
    $ws.Rows.Item(321).Insert()
}

# --- Row 321: CCSS_AMARIA_ALUPAR (new) ---
$ws.Range("A321").Value = "CCSS_AMARIA_ALUPAR"
$ws.Range("B321").Value = "-"
$ws.Range("C321").Value = "-"
$ws.Range("D321").Value = ""
$ws.Range("E321").Value = "Generador Sincrónico"
$ws.Range("F321").Value = 0
$ws.Range("G321").Value = "LF5"
$ws.Range("H321").Value = "00-Norte Grande"
$ws.Range("I321").Value = ""
$ws.Range("J321").Value = ""
$ws.Range("K321").Value = ""
$ws.Range("L321").Value = "PV"
$ws.Range("M321").Value = 360
$ws.Range("N321").Value = 0
$ws.Range("O321").Value = 0
$ws.Range("P321").Value = 15
$ws.Range("Q321").Value = 15
$ws.Range("R321").Value = 360
Set-TextValue "S321" "1"
Set-TextValue "T321" "1"

# --- Row 322: CCSS_AMARIA_TRANSELEC (new) ---
$ws.Range("A322").Value = "CCSS_AMARIA_TRANSELEC"
$ws.Range("B322").Value = "-"
$ws.Range("C322").Value = "-"
$ws.Range("D322").Value = ""
$ws.Range("E322").Value = "Generador Sincrónico"
$ws.Range("F322").Value = 0
$ws.Range("G322").Value = "LF7"
$ws.Range("H322").Value = "00-Norte Grande"
$ws.Range("I322").Value = ""
$ws.Range("J322").Value = ""
$ws.Range("K322").Value = ""
$ws.Range("L322").Value = "PV"
$ws.Range("M322").Value = 92
$ws.Range("N322").Value = 0
$ws.Range("O322").Value = 0
$ws.Range("P322").Value = 15
$ws.Range("Q322").Value = 15
$ws.Range("R322").Value = 92
Set-TextValue "S322" "1"
Set-TextValue "T322" "1"

# --- Row 323: CCSS_LIKANANTAI_TRANSELEC (new) ---
$ws.Range("A323").Value = "CCSS_LIKANANTAI_TRANSELEC"
$ws.Range("B323").Value = "-"
$ws.Range("C323").Value = "-"
$ws.Range("D323").Value = ""
$ws.Range("E323").Value = "Generador Sincrónico"
$ws.Range("F323").Value = 0
$ws.Range("G323").Value = "LF6"
$ws.Range("H323").Value = "00-Norte Grande"
$ws.Range("I323").Value = ""
$ws.Range("J323").Value = ""
$ws.Range("K323").Value = ""
$ws.Range("L323").Value = "PV"
$ws.Range("M323").Value = 92
$ws.Range("N323").Value = 0
$ws.Range("O323").Value = 0
$ws.Range("P323").Value = 15
$ws.Range("Q323").Value = 15
$ws.Range("R323").Value = 92
Set-TextValue "S323" "1"
Set-TextValue "T323" "1"

# --- Row 324: CCSS_ILLAPA_ALUPAR (new) ---
$ws.Range("A324").Value = "CCSS_ILLAPA_ALUPAR"
$ws.Range("B324").Value = "-"
$ws.Range("C324").Value = "-"
$ws.Range("D324").Value = ""
$ws.Range("E324").Value = "Generador Sincrónico"
$ws.Range("F324").Value = 0
$ws.Range("G324").Value = "LF4"
$ws.Range("H324").Value = "01-Atacama"
$ws.Range("I324").Value = ""
$ws.Range("J324").Value = ""
$ws.Range("K324").Value = ""
$ws.Range("L324").Value = "PV"
$ws.Range("M324").Value = 270
$ws.Range("N324").Value = 0
$ws.Range("O324").Value = 0
$ws.Range("P324").Value = 15
$ws.Range("Q324").Value = 15
$ws.Range("R324").Value = 270
Set-TextValue "S324" "1"
Set-TextValue "T324" "1"

# --- Row 325: CCSS_TOCOPILLA_ENGIE (new) ---
$ws.Range("A325").Value = "CCSS_TOCOPILLA_ENGIE"
$ws.Range("B325").Value = "-"
$ws.Range("C325").Value = "-"
$ws.Range("D325").Value = ""
$ws.Range("E325").Value = "Generador Sincrónico"
Set-TextValue "F325" "1"
$ws.Range("G325").Value = "LF4"
$ws.Range("H325").Value = "00-Norte Grande"
$ws.Range("I325").Value = ""
$ws.Range("J325").Value = ""
$ws.Range("K325").Value = ""
$ws.Range("L325").Value = "PV"
$ws.Range("M325").Value = 147.1
$ws.Range("N325").Value = 0
$ws.Range("O325").Value = 0
$ws.Range("P325").Value = 15
$ws.Range("Q325").Value = 13.8
$ws.Range("R325").Value = 147.1
Set-TextValue "S325" "1"
Set-TextValue "T325" "1"

# --- Row 326: STAT_PAZUCAR (shifted from old row 321) -> J changes -0 to 0 ---
$ws.Range("J326").Value = 0

# --- Row 327: STAT_POLPAICO (shifted from old row 322) -> unchanged ---

# --- Row 328: ILLAPA_GFM (shifted from old row 323) -> F,J,K reset to 0 ---
$ws.Range("F328").Value = 0
$ws.Range("J328").Value = 0
$ws.Range("K328").Value = 0

# --- Row 329: LIKANA_GFM (shifted from old row 324) -> F,J,K reset to 0 ---
$ws.Range("F329").Value = 0
$ws.Range("J329").Value = 0
$ws.Range("K329").Value = 0

# --- Row 330: NCHUQUI_GFM (shifted from old row 325) -> F,J,K reset to 0 ---
$ws.Range("F330").Value = 0
$ws.Range("J330").Value = 0
$ws.Range("K330").Value = 0

# --- Row 331: AMARIA_GFM (shifted from old row 326) -> F,J,K reset to 0 ---
$ws.Range("F331").Value = 0
$ws.Range("J331").Value = 0
$ws.Range("K331").Value = 0

# --- Row 332: STAT_CNAVIA (shifted from old row 327) -> unchanged ---
